# Append the new daily portfolio row (row 59) to the worksheet, matching the
# values added on 2025-10-13: Date, SUZLON.NS, TATAMOTORS.NS, ETERNAL.NS.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates stored as plain text (e.g. "2025-10-12" in row 58),
# not as date serials, so force text formatting before assigning the value
# to stop Excel from auto-converting the string into a date number. Reset
# the style back to Normal afterwards so the cell keeps the same (default)
# formatting as the rest of the column.
$ws.Cells.Item(59, 1).NumberFormat = "@"
$ws.Cells.Item(59, 1).Value = "2025-10-13"
$ws.Cells.Item(59, 1).Style = "Normal"

$ws.Cells.Item(59, 2).Value = 54.34999847412109
$ws.Cells.Item(59, 3).Value = 660.75
$ws.Cells.Item(59, 4).Value = 348.3500061035156
